# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.440.53"
$ws.Range("E2").Value = "  +1.09%  "
# Row 3
$ws.Range("D3").Value = "3.809.13"
$ws.Range("E3").Value = "  +0.21%  "
# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.26%  "
# Row 5
$ws.Range("D5").Value = "'608.13"
$ws.Range("E5").Value = "  +0.92%  "
# Row 6
$ws.Range("D6").Value = "'164.22"
$ws.Range("E6").Value = "  -0.91%  "
# Row 7
$ws.Range("D7").Value = "3.807.46"
$ws.Range("E7").Value = "  +0.33%  "
# Row 8
$ws.Range("E8").Value = "  -0.15%  "
# Row 9
$ws.Range("E9").Value = "  -0.12%  "
# Row 10
$ws.Range("E10").Value = "  +0.21%  "
# Row 11
$ws.Range("D11").Value = "'7.00"
$ws.Range("E11").Value = "  +11.52%  "
# Row 12
$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = "  +0.15%  "
# Row 13
$ws.Range("E13").Value = "  -0.60%  "
# Row 14
$ws.Range("D14").Value = "'35.17"
$ws.Range("E14").Value = "  -1.69%  "
# Row 15
$ws.Range("D15").Value = "4.446.09"
$ws.Range("E15").Value = "  +0.02%  "
# Row 16
$ws.Range("D16").Value = "3.822.18"
$ws.Range("E16").Value = "  +0.28%  "
# Row 17
$ws.Range("D17").Value = "68.434.24"
$ws.Range("E17").Value = "  +1.05%  "
# Row 18
$ws.Range("D18").Value = "'18.13"
$ws.Range("E18").Value = "  -1.82%  "
# Row 19
$ws.Range("E19").Value = "  +1.90%  "
# Row 20
$ws.Range("E20").Value = "  +0.46%  "
# Row 21
$ws.Range("D21").Value = "'463.48"
$ws.Range("E21").Value = "  +0.31%  "
# Row 22
$ws.Range("E22").Value = "  -2.21%  "
# Row 23
$ws.Range("D23").Value = "'0.700"
$ws.Range("E23").Value = "  +0.13%  "
# Row 24
$ws.Range("E24").Value = "  +1.12%  "
# Row 25
$ws.Range("D25").Value = "'83.64"
$ws.Range("E25").Value = "  +0.56%  "
# Row 26
$ws.Range("E26").Value = "  -0.54%  "
# Row 27
$ws.Range("D27").Value = "'2.11"
$ws.Range("E27").Value = "  -0.05%  "
# Row 28
$ws.Range("E28").Value = "  -0.20%  "
# Row 29
$ws.Range("E29").Value = "  +0.12%  "
# Row 30
$ws.Range("D30").Value = "3.957.15"
$ws.Range("E30").Value = "  +0.13%  "
# Row 31
$ws.Range("E31").Value = "  -5.30%  "
# Row 32
$ws.Range("D32").Value = "'2.23"
$ws.Range("E32").Value = "  +0.58%  "
# Row 33
$ws.Range("E33").Value = "  -1.01%  "
# Row 34
$ws.Range("D34").Value = "'29.17"
$ws.Range("E34").Value = "  -0.73%  "
# Row 35
$ws.Range("E35").Value = "  -0.13%  "
# Row 36
$ws.Range("D36").Value = "'9.04"
$ws.Range("E36").Value = "  -0.28%  "
# Row 37
$ws.Range("E37").Value = "  +1.44%  "
# Row 38
$ws.Range("E38").Value = "  +8.20%  "
# Row 39
$ws.Range("D39").Value = "'5.90"
$ws.Range("E39").Value = "  +1.82%  "
# Row 40
$ws.Range("D40").Value = "'3.23"
$ws.Range("E40").Value = "  -0.14%  "
# Row 41
$ws.Range("D41").Value = "'0.981"
$ws.Range("E41").Value = "  -1.55%  "
# Row 42
$ws.Range("E42").Value = "  -0.02%  "
# Row 43
$ws.Range("E43").Value = "  +0.01%  "
# Row 44
$ws.Range("E44").Value = "  -0.46%  "
# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'47.11"
$ws.Range("E45").Value = "  -1.12%  "
# Row 46
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'153.10"
$ws.Range("E46").Value = "  +1.55%  "
# Row 47
$ws.Range("D47").Value = "'43.13"
$ws.Range("E47").Value = "  -3.21%  "
# Row 48
$ws.Range("E48").Value = "  -0.39%  "
# Row 49
$ws.Range("E49").Value = "  +0.78%  "
# Row 50
$ws.Range("E50").Value = "  +0.98%  "
# Row 51
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'26.26"
$ws.Range("E51").Value = "  -7.52%  "
